$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D column cluster labels for rows 5, 6, 11, 12
# (text values are unchanged: D5/D11 = MuSCs, D6/D12 = Neutrophils;
#  this keeps the shared-string table reorder consistent when Excel resaves)
$ws.Range("D5").Value = "MuSCs"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("D11").Value = "MuSCs"
$ws.Range("D12").Value = "Neutrophils"

# Update recomputed TPM-derived expression values (rows 2-13)
# Row 2
$ws.Range("G2").Value = 0.3750503333333333
$ws.Range("H2").Value = 1.125151
$ws.Range("I2").Value = 0.7477030322765
$ws.Range("J2").Value = 0.7477030322765
$ws.Range("M2").Value = 156.53184
$ws.Range("N2").Value = 469.59552
$ws.Range("O2").Value = 0.1403721039197297
$ws.Range("P2").Value = 0.1403721039197297
$ws.Range("Q2").Value = 58.70731876927999
$ws.Range("R2").Value = 528.36586892352
$ws.Range("S2").Value = 0.1049566477478139
$ws.Range("T2").Value = 0.1049566477478139

# Row 3
$ws.Range("G3").Value = 0.3750503333333333
$ws.Range("H3").Value = 1.125151
$ws.Range("I3").Value = 0.7477030322765
$ws.Range("J3").Value = 0.7477030322765
$ws.Range("O3").Value = 0.3490212885850074
$ws.Range("P3").Value = 0.3490212885850074
$ws.Range("Q3").Value = 145.9699147769555
$ws.Range("R3").Value = 1313.7292329926
$ws.Range("S3").Value = 0.2609642758040614
$ws.Range("T3").Value = 0.2609642758040614

# Row 4
$ws.Range("G4").Value = 0.3750503333333333
$ws.Range("H4").Value = 1.125151
$ws.Range("I4").Value = 0.7477030322765
$ws.Range("J4").Value = 0.7477030322765
$ws.Range("M4").Value = 169.4499613333333
$ws.Range("N4").Value = 508.349884
$ws.Range("O4").Value = 0.1519566088373896
$ws.Range("P4").Value = 0.1519566088373896
$ws.Range("Q4").Value = 63.55226448138711
$ws.Range("R4").Value = 571.970380332484
$ws.Range("S4").Value = 0.1136184172021702
$ws.Range("T4").Value = 0.1136184172021702

# Row 5
$ws.Range("G5").Value = 0.3750503333333333
$ws.Range("H5").Value = 1.125151
$ws.Range("I5").Value = 0.7477030322765
$ws.Range("J5").Value = 0.7477030322765
$ws.Range("M5").Value = 169.0002543333333
$ws.Range("N5").Value = 507.000763
$ws.Range("O5").Value = 0.1515533278324679
$ws.Range("P5").Value = 0.1515533278324679
$ws.Range("Q5").Value = 63.38360172113478
$ws.Range("R5").Value = 570.452415490213
$ws.Range("S5").Value = 0.1133168827719308
$ws.Range("T5").Value = 0.1133168827719308

# Row 6
$ws.Range("G6").Value = 0.3750503333333333
$ws.Range("H6").Value = 1.125151
$ws.Range("I6").Value = 0.7477030322765
$ws.Range("J6").Value = 0.7477030322765
$ws.Range("M6").Value = 60.30985666666667
$ws.Range("N6").Value = 180.92957
$ws.Range("O6").Value = 0.05408370250677011
$ws.Range("P6").Value = 0.05408370250677011
$ws.Range("Q6").Value = 22.61923184611889
$ws.Range("R6").Value = 203.57308661507
$ws.Range("S6").Value = 0.04043854836105216
$ws.Range("T6").Value = 0.04043854836105216

# Row 7
$ws.Range("G7").Value = 0.3750503333333333
$ws.Range("H7").Value = 1.125151
$ws.Range("I7").Value = 0.7477030322765
$ws.Range("J7").Value = 0.7477030322765
$ws.Range("M7").Value = 170.6279296666667
$ws.Range("N7").Value = 511.883789
$ws.Range("O7").Value = 0.1530129683186351
$ws.Range("P7").Value = 0.1530129683186351
$ws.Range("Q7").Value = 63.99406189745989
$ws.Range("R7").Value = 575.946557077139
$ws.Range("S7").Value = 0.1144082603894715
$ws.Range("T7").Value = 0.1144082603894715

# Row 8
$ws.Range("G8").Value = 0.126553
$ws.Range("H8").Value = 0.379659
$ws.Range("I8").Value = 0.2522969677235
$ws.Range("J8").Value = 0.2522969677235
$ws.Range("M8").Value = 156.53184
$ws.Range("N8").Value = 469.59552
$ws.Range("O8").Value = 0.1403721039197297
$ws.Range("P8").Value = 0.1403721039197297
$ws.Range("Q8").Value = 19.80957394752
$ws.Range("R8").Value = 178.28616552768
$ws.Range("S8").Value = 0.03541545617191583
$ws.Range("T8").Value = 0.03541545617191583

# Row 9
$ws.Range("G9").Value = 0.126553
$ws.Range("H9").Value = 0.379659
$ws.Range("I9").Value = 0.2522969677235
$ws.Range("J9").Value = 0.2522969677235
$ws.Range("O9").Value = 0.3490212885850074
$ws.Range("P9").Value = 0.3490212885850074
$ws.Range("Q9").Value = 49.25453727926666
$ws.Range("R9").Value = 443.2908355134
$ws.Range("S9").Value = 0.08805701278094599
$ws.Range("T9").Value = 0.08805701278094599

# Row 10
$ws.Range("G10").Value = 0.126553
$ws.Range("H10").Value = 0.379659
$ws.Range("I10").Value = 0.2522969677235
$ws.Range("J10").Value = 0.2522969677235
$ws.Range("M10").Value = 169.4499613333333
$ws.Range("N10").Value = 508.349884
$ws.Range("O10").Value = 0.1519566088373896
$ws.Range("P10").Value = 0.1519566088373896
$ws.Range("Q10").Value = 21.44440095661733
$ws.Range("R10").Value = 192.999608609556
$ws.Range("S10").Value = 0.03833819163521939
$ws.Range("T10").Value = 0.03833819163521939

# Row 11
$ws.Range("G11").Value = 0.126553
$ws.Range("H11").Value = 0.379659
$ws.Range("I11").Value = 0.2522969677235
$ws.Range("J11").Value = 0.2522969677235
$ws.Range("M11").Value = 169.0002543333333
$ws.Range("N11").Value = 507.000763
$ws.Range("O11").Value = 0.1515533278324679
$ws.Range("P11").Value = 0.1515533278324679
$ws.Range("Q11").Value = 21.38748918664633
$ws.Range("R11").Value = 192.487402679817
$ws.Range("S11").Value = 0.03823644506053718
$ws.Range("T11").Value = 0.03823644506053718

# Row 12
$ws.Range("G12").Value = 0.126553
$ws.Range("H12").Value = 0.379659
$ws.Range("I12").Value = 0.2522969677235
$ws.Range("J12").Value = 0.2522969677235
$ws.Range("M12").Value = 60.30985666666667
$ws.Range("N12").Value = 180.92957
$ws.Range("O12").Value = 0.05408370250677011
$ws.Range("P12").Value = 0.05408370250677011
$ws.Range("Q12").Value = 7.632393290736667
$ws.Range("R12").Value = 68.69153961663001
$ws.Range("S12").Value = 0.01364515414571795
$ws.Range("T12").Value = 0.01364515414571795

# Row 13
$ws.Range("G13").Value = 0.126553
$ws.Range("H13").Value = 0.379659
$ws.Range("I13").Value = 0.2522969677235
$ws.Range("J13").Value = 0.2522969677235
$ws.Range("M13").Value = 170.6279296666667
$ws.Range("N13").Value = 511.883789
$ws.Range("O13").Value = 0.1530129683186351
$ws.Range("P13").Value = 0.1530129683186351
$ws.Range("Q13").Value = 21.59347638310567
$ws.Range("R13").Value = 194.341287447951
$ws.Range("S13").Value = 0.03860470792916362
$ws.Range("T13").Value = 0.03860470792916362
